$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value without Excel coercing it to a number,
# and without leaving a permanent number-format style on the cell.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "245.67"
Set-TextValue $ws "D4" "5.406"
Set-TextValue $ws "D5" "0.05763"
Set-TextValue $ws "D6" "3.432"
Set-TextValue $ws "D7" "6.331"
Set-TextValue $ws "D8" "0.8100"
Set-TextValue $ws "D9" "0.8945"
Set-TextValue $ws "D10" "0.1447"
Set-TextValue $ws "D11" "0.07419"
Set-TextValue $ws "D12" "0.03127"
$ws.Range("E12").Value = "11BitrueCoinBTRBestin24h"
Set-TextValue $ws "D13" "0.02988"
Set-TextValue $ws "D14" "0.09419"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D15" "3.934"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D16" "0.001580"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D17" "0.04784"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D18" "0.0005851"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D19" "0.006337"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D20" "0.004070"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D21" "0.0009971"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws "D22" "0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D23" "3.733"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D24" "2.197"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws "D25" "0.3273"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws "D26" "0.1300"
$ws.Range("E26").Value = "25ProBitTokenPROB"
Set-TextValue $ws "D27" "0.0004651"
Set-TextValue $ws "D40" "0.03904"
Set-TextValue $ws "D41" "0.006789"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue $ws "D42" "0.1073"
Set-TextValue $ws "D43" "0.002431"
Set-TextValue $ws "D44" "0.006817"
Set-TextValue $ws "D45" "0.00005650"
Set-TextValue $ws "D47" "0.3801"
Set-TextValue $ws "D48" "0.1629"
